$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-31, columns C:H (ax, ay, az, gx, gy, gz)
$data = @(
    @(-1.260086536407471, -6.693731784820557, -10.91690349578857, 0.09952953143403931, -0.1946007158231257, -1.235400125247621),
    @(-0.3250093758106231, -9.289229393005373, -12.17654037475586, 0.7596392552240191, -0.4710539486430101, -1.249458708894362),
    @(-0.197640836238861, -8.023595809936523, -9.897425651550291, 0.9691518504138372, -1.144679769326803, -1.013632127177822),
    @(8.160659790039062, -11.02071857452393, -0.9127321243286132, 0.5637087537607993, -0.8698977423370384, -0.6591360478773013),
    @(0.2106014490127563, -8.135941505432129, -7.824718952178955, -1.529567692257933, -0.1701191907423907, 2.089324393540335),
    @(3.468891382217407, -25.93060493469238, 3.170974254608154, -1.60695329718632, 0.9596939690479577, 2.503135678965),
    @(-34.4773063659668, -0.4352650642395019, -17.56653213500977, -4.042883170033875, 4.578179887675418, 0.4566991565697582),
    @(-14.77978038787842, -11.52226448059082, -11.90310859680176, -0.6153825116813314, 0.4450490321587646, -1.158701128369084),
    @(24.58778762817383, -1.438390731811523, 15.80244827270508, 5.463495351852635, -3.977313330414098, -3.230007587222884),
    @(3.9564049243927, -11.81492042541504, 11.41962051391602, 4.66184952062202, 3.56197387800312, -1.986260027786983),
    @(0.3332785367965698, 0.8126574158668518, -8.991792678833008, -0.716348304661115, 4.782481053553503, 2.407147008686489),
    @(26.99358177185058, -74.40792083740234, 2.843526840209961, -4.534485773209054, 1.10340199667379, 3.735504924704164),
    @(47.74433135986328, 8.593025207519531, -17.67363357543945, -5.989971771152742, 4.972828354857392, 1.774531464511098),
    @(0.0659644603729248, -10.85460472106934, -14.575852394104, -4.013937945759595, 5.69700783545815, -1.538702366548966),
    @(29.60918045043945, -78.06973266601562, 16.24357986450195, 2.29124034763484, 0.5539356289653035, -3.10506863769041),
    @(-4.278045177459717, 2.075074911117554, -9.626623153686523, 4.105864645144252, -2.789147576060868, -3.677538324933537),
    @(16.6558837890625, -14.57428646087646, 21.4177131652832, 0.2592948106450716, -3.157548316028139, -0.8597697675774869),
    @(-54.05896377563477, -32.270751953125, -11.78435802459717, -6.645317213251005, -0.560946681083867, 5.571864681506409),
    @(1.290312767028809, -2.825088739395142, -12.40687561035156, -1.937935614804442, -0.1515454994429453, -1.15992745128265),
    @(-6.366491794586182, -13.87612628936768, -13.16234016418457, 2.910339694743533, -5.836961925576653, -4.45155078828873),
    @(-4.204625129699707, -2.39580774307251, -2.694075107574463, 3.819924138400946, -3.973247524795179, -0.4456482636272369),
    @(3.387409210205078, -0.1273889541625976, -3.40644645690918, 1.989335611325833, -0.2693817703002244, 0.09849660987155473),
    @(-30.62431526184082, -35.95574188232422, -0.9927592277526855, -3.755865255627084, 0.09460148078591657, 3.946894074251968),
    @(31.11456108093262, 10.69985771179199, -13.1208086013794, -2.889570736010117, 9.743136939652585, 0.5475007984615461),
    @(-20.68763160705566, -28.60056304931641, -20.36888885498047, -1.842574656009669, 4.510396480560456, -1.828126966953318),
    @(27.19162559509277, -76.43211364746094, -1.331388473510742, 2.88307266279109, -8.930636825911506, -1.004032475139075),
    @(2.005526781082153, 2.733113288879395, -6.29358959197998, 4.052329929460007, -3.608018367662238, -3.160855819325922),
    @(-72.39521789550781, -77.2735824584961, 66.16534423828125, -1.29231422538059, 1.658297274637648, -0.279984981641821),
    @(-17.93119621276855, 35.98200607299805, 4.464548587799072, -4.269718651377847, -2.122582794056059, -2.14775515915062),
    @(-0.0441454052925109, -11.34683990478516, -10.79495239257812, -2.918748300009885, 1.177203515254081, -6.725167151984817)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = 3 + $j  # column C = 3
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}

# Append new rows 22-31 with timestamp (A) and label (B)
$timestamps = @(2000, 2100, 2200, 2300, 2400, 2500, 2600, 2700, 2800, 2900)
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = 22 + $i
    $ws.Cells.Item($row, 1).Value = $timestamps[$i]
    $ws.Cells.Item($row, 2).Value = "walkingToRunning"
}